$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the credential values stored in A2/B2
$ws.Range("A2").Value = "ashish dhaundi"
$ws.Range("B2").Value = "Nimda@1234"

# Move the active selection to B2 (was previously E9)
$ws.Range("B2").Select()
